# HKStock: add function readIndexData -- append three more trading days
# (2016-08-09, 2016-08-10, 2016-08-11) to the HSF index history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# New rows to append right after the existing last row (264).
$newRows = @(
    @{ Row = 265; Idx = 263; Date = "2016-08-09"; Open = 29248.38; Low = 29120.68; High = 29294.86; Close = 29230.62 },
    @{ Row = 266; Idx = 264; Date = "2016-08-10"; Open = 29341.55; Low = 29212.38; High = 29479.93; Close = 29317.3 },
    @{ Row = 267; Idx = 265; Date = "2016-08-11"; Open = 29179.59; Low = 29145.45; High = 29802.42; Close = 29629.76 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $srcRow = $row - 1

    # Clone the whole row's formatting from the row above it so the new
    # row's styling (bold/bordered index column, etc.) matches the rest
    # of the table, without introducing any new cell-style entries.
    $ws.Range("A" + $srcRow + ":G" + $srcRow).Copy()
    $ws.Range("A" + $row + ":G" + $row).PasteSpecial($xlPasteFormats)

    $ws.Cells.Item($row, 1).Value = $r.Idx
    $ws.Cells.Item($row, 2).Value = "HSF"

    # Write the trade date as a literal text formula, then flatten it to a
    # static value via copy/paste-values. This keeps the cell a plain
    # shared-string (matching the other tradeDate cells) instead of Excel
    # auto-converting the "YYYY-MM-DD"-looking text into a date serial.
    $c = $ws.Cells.Item($row, 3)
    $c.Formula = '="' + $r.Date + '"'
    $c.Copy()
    $c.PasteSpecial($xlPasteValues)

    $ws.Cells.Item($row, 4).Value = $r.Open
    $ws.Cells.Item($row, 5).Value = $r.Low
    $ws.Cells.Item($row, 6).Value = $r.High
    $ws.Cells.Item($row, 7).Value = $r.Close
}
